$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 472.1111
$ws.Range("I41").Value = 313.1
$ws.Range("J41").Value = 533.2692
$ws.Range("K41").Value = 313.1
$ws.Range("L41").Value = 533.2692
$ws.Range("M41").Value = 126.9
$ws.Range("N41").Value = -1413.2692
# Row 46
$ws.Range("H46").Value = 1907.4615
$ws.Range("I46").Value = 1017
$ws.Range("J46").Value = 1981.6666
$ws.Range("K46").Value = 3051
$ws.Range("L46").Value = 5944.9998
$ws.Range("M46").Value = -2932
$ws.Range("N46").Value = -6182.9998
# Row 60
$ws.Range("H60").Value = 1907.4615
$ws.Range("I60").Value = 1017
$ws.Range("J60").Value = 1981.6666
$ws.Range("K60").Value = 3051
$ws.Range("L60").Value = 5944.9998
$ws.Range("M60").Value = -2567
$ws.Range("N60").Value = -6912.9998
# Row 111
$ws.Range("H111").Value = 394.875
$ws.Range("I111").Value = 359.75
$ws.Range("J111").Value = 430
$ws.Range("K111").Value = 1079.25
$ws.Range("L111").Value = 1290
$ws.Range("M111").Value = 1987.75
$ws.Range("N111").Value = -7424
# Row 125
$ws.Range("H125").Value = 1342.5883
$ws.Range("I125").Value = 1022.1818
$ws.Range("K125").Value = 9199.6362
$ws.Range("M125").Value = -6739.6362
# Row 141
$ws.Range("H141").Value = 5209.5835
$ws.Range("I141").Value = 6195.7144
$ws.Range("J141").Value = 3829
$ws.Range("K141").Value = 18587.1432
$ws.Range("L141").Value = 11487
$ws.Range("M141").Value = -13407.1432
$ws.Range("N141").Value = -21847

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 964.8889
$ws.Range("I2").Value = 814.6
$ws.Range("J2").Value = 1394.2858
$ws.Range("K2").Value = 814.6
$ws.Range("L2").Value = 1394.2858
$ws.Range("M2").Value = -701.6
$ws.Range("N2").Value = -1620.2858
# Row 7
$ws.Range("H7").Value = 39653
$ws.Range("J7").Value = 39653
$ws.Range("L7").Value = 39653
$ws.Range("N7").Value = -39881
# Row 32
$ws.Range("H32").Value = 21445.53
$ws.Range("I32").Value = 3685.6575
$ws.Range("K32").Value = 3685.6575
$ws.Range("M32").Value = -3398.6575
# Row 34
$ws.Range("H34").Value = 32500
$ws.Range("I34").Value = 20000
$ws.Range("K34").Value = 20000
$ws.Range("M34").Value = -19729
# Row 45
$ws.Range("H45").Value = 1656
$ws.Range("I45").Value = 1604.3478
$ws.Range("J45").Value = 2250
$ws.Range("K45").Value = 1604.3478
$ws.Range("L45").Value = 2250
$ws.Range("M45").Value = -1227.3478
$ws.Range("N45").Value = -3004
# Row 116
$ws.Range("H116").Value = 964.8889
$ws.Range("I116").Value = 814.6
$ws.Range("J116").Value = 1394.2858
$ws.Range("K116").Value = 814.6
$ws.Range("L116").Value = 1394.2858
$ws.Range("M116").Value = 1479.4
$ws.Range("N116").Value = -5982.2858
# Row 132
$ws.Range("H132").Value = 2125.2856
$ws.Range("I132").Value = 1806.7778
$ws.Range("J132").Value = 3558.5715
$ws.Range("K132").Value = 5420.3334
$ws.Range("L132").Value = 10675.7145
$ws.Range("M132").Value = -2890.3334
$ws.Range("N132").Value = -15735.7145

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 964.8889
$ws.Range("I3").Value = 814.6
$ws.Range("J3").Value = 1394.2858
$ws.Range("K3").Value = 814.6
$ws.Range("L3").Value = 1394.2858
$ws.Range("M3").Value = -700.6
$ws.Range("N3").Value = -1622.2858
# Row 55
$ws.Range("H55").Value = 38500
$ws.Range("J55").Value = 38500
$ws.Range("L55").Value = 38500
$ws.Range("N55").Value = -39046
# Row 94
$ws.Range("H94").Value = 734.65625
$ws.Range("I94").Value = 704.3043
$ws.Range("J94").Value = 812.2222
$ws.Range("K94").Value = 704.3043
$ws.Range("L94").Value = 812.2222
$ws.Range("M94").Value = -253.3043
$ws.Range("N94").Value = -1714.2222
# Row 107
$ws.Range("H107").Value = 1136.4286
$ws.Range("I107").Value = 1136.4286
$ws.Range("K107").Value = 1136.4286
$ws.Range("M107").Value = 783.5714
# Row 139
$ws.Range("H139").Value = 30280
$ws.Range("J139").Value = 30280
$ws.Range("L139").Value = 30280
$ws.Range("N139").Value = -40560

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 35
$ws.Range("H35").Value = 14955
$ws.Range("I35").Value = 1193.75
$ws.Range("J35").Value = 70000
$ws.Range("K35").Value = 1193.75
$ws.Range("L35").Value = 70000
$ws.Range("M35").Value = -899.75
$ws.Range("N35").Value = -70588
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 121
$ws.Range("H121").Value = 50000
$ws.Range("J121").Value = 50000
$ws.Range("L121").Value = 50000
$ws.Range("N121").Value = -52620
# Row 122
$ws.Range("H122").Value = 3087.6667
$ws.Range("I122").Value = 2882.4285
$ws.Range("J122").Value = 3375
$ws.Range("K122").Value = 8647.2855
$ws.Range("L122").Value = 10125
$ws.Range("M122").Value = -6197.2855
$ws.Range("N122").Value = -15025
# Row 132
$ws.Range("H132").Value = 1398.6531
$ws.Range("I132").Value = 1094.55
$ws.Range("J132").Value = 2750.2222
$ws.Range("K132").Value = 3283.65
$ws.Range("L132").Value = 8250.6666
$ws.Range("M132").Value = -753.6499999999996
$ws.Range("N132").Value = -13310.6666

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 901.5
$ws.Range("I22").Value = 403
$ws.Range("J22").Value = 1400
$ws.Range("K22").Value = 1209
$ws.Range("L22").Value = 4200
$ws.Range("M22").Value = -1040
$ws.Range("N22").Value = -4538
# Row 27
$ws.Range("H27").Value = 901.5
$ws.Range("I27").Value = 403
$ws.Range("J27").Value = 1400
$ws.Range("K27").Value = 1209
$ws.Range("L27").Value = 4200
$ws.Range("M27").Value = -1107
$ws.Range("N27").Value = -4404
# Row 107
$ws.Range("H107").Value = 5797.8945
$ws.Range("I107").Value = 17161
$ws.Range("J107").Value = 553.38464
$ws.Range("K107").Value = 51483
$ws.Range("L107").Value = 1660.15392
$ws.Range("M107").Value = -49563
$ws.Range("N107").Value = -5500.15392
# Row 113
$ws.Range("H113").Value = 2135.111
$ws.Range("J113").Value = 491.66666
$ws.Range("L113").Value = 1474.99998
$ws.Range("N113").Value = -5814.999980000001
# Row 122
$ws.Range("H122").Value = 1438.5588
$ws.Range("I122").Value = 1385.6154
$ws.Range("J122").Value = 1471.3334
$ws.Range("K122").Value = 12470.5386
$ws.Range("L122").Value = 13242.0006
$ws.Range("M122").Value = -10020.5386
$ws.Range("N122").Value = -18142.0006

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4541.591
$ws.Range("I70").Value = 4472.5713
$ws.Range("J70").Value = 4662.375
$ws.Range("K70").Value = 4472.5713
$ws.Range("L70").Value = 4662.375
$ws.Range("M70").Value = -4202.5713
$ws.Range("N70").Value = -5202.375
# Row 73
$ws.Range("H73").Value = 4541.591
$ws.Range("I73").Value = 4472.5713
$ws.Range("J73").Value = 4662.375
$ws.Range("K73").Value = 4472.5713
$ws.Range("L73").Value = 4662.375
$ws.Range("M73").Value = -3536.5713
$ws.Range("N73").Value = -6534.375
# Row 102
$ws.Range("H102").Value = 1083.7
$ws.Range("I102").Value = 1083.7
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1083.7
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 538.3
$ws.Range("N102").ClearContents()
# Row 132
$ws.Range("H132").Value = 3001.077
$ws.Range("I132").Value = 2669.8125
$ws.Range("J132").Value = 3531.1
$ws.Range("K132").Value = 8009.4375
$ws.Range("L132").Value = 10593.3
$ws.Range("M132").Value = -5479.4375
$ws.Range("N132").Value = -15653.3

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1908.1428
$ws.Range("I7").Value = 1437.375
$ws.Range("J7").Value = 2535.8333
$ws.Range("K7").Value = 1437.375
$ws.Range("L7").Value = 2535.8333
$ws.Range("M7").Value = -1325.375
$ws.Range("N7").Value = -2759.8333
# Row 32
$ws.Range("H32").Value = 10328.625
$ws.Range("I32").Value = 2199.6667
$ws.Range("J32").Value = 15206
$ws.Range("K32").Value = 2199.6667
$ws.Range("L32").Value = 15206
$ws.Range("M32").Value = -1882.6667
$ws.Range("N32").Value = -15840
# Row 40
$ws.Range("H40").Value = 2425.818
$ws.Range("I40").Value = 2309.3333
$ws.Range("J40").Value = 2950
$ws.Range("K40").Value = 2309.3333
$ws.Range("L40").Value = 2950
$ws.Range("M40").Value = -2173.3333
$ws.Range("N40").Value = -3222
# Row 46
$ws.Range("H46").Value = 650933.8
$ws.Range("I46").Value = 600
$ws.Range("J46").Value = 781000.6
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 781000.6
$ws.Range("M46").Value = -412
$ws.Range("N46").Value = -781376.6
# Row 61
$ws.Range("H61").Value = 2078.5715
$ws.Range("I61").Value = 1504.2858
$ws.Range("J61").Value = 2652.8572
$ws.Range("K61").Value = 1504.2858
$ws.Range("L61").Value = 2652.8572
$ws.Range("M61").Value = -1302.2858
$ws.Range("N61").Value = -3056.8572
# Row 113
$ws.Range("H113").Value = 2078.5715
$ws.Range("I113").Value = 1504.2858
$ws.Range("J113").Value = 2652.8572
$ws.Range("K113").Value = 1504.2858
$ws.Range("L113").Value = 2652.8572
$ws.Range("M113").Value = 665.7141999999999
$ws.Range("N113").Value = -6992.8572
# Row 122
$ws.Range("H122").Value = 2966.3901
$ws.Range("I122").Value = 2853.5881
$ws.Range("J122").Value = 3514.2856
$ws.Range("K122").Value = 8560.764299999999
$ws.Range("L122").Value = 10542.8568
$ws.Range("M122").Value = -6110.764299999999
$ws.Range("N122").Value = -15442.8568
# Row 126
$ws.Range("H126").Value = 1908.1428
$ws.Range("I126").Value = 1437.375
$ws.Range("J126").Value = 2535.8333
$ws.Range("K126").Value = 4312.125
$ws.Range("L126").Value = 7607.499899999999
$ws.Range("M126").Value = -1842.125
$ws.Range("N126").Value = -12547.4999

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 82
$ws.Range("H82").Value = 20000
$ws.Range("J82").Value = 20000
$ws.Range("L82").Value = 20000
$ws.Range("N82").Value = -20766
# Row 85
$ws.Range("H85").Value = 20000
$ws.Range("J85").Value = 20000
$ws.Range("L85").Value = 20000
$ws.Range("N85").Value = -22652
# Row 126
$ws.Range("H126").Value = 250573.08
$ws.Range("I126").Value = 286132.38
$ws.Range("K126").Value = 858397.14
$ws.Range("M126").Value = -855927.14
# Row 129
$ws.Range("H129").Value = 40427.625
$ws.Range("J129").Value = 40427.625
$ws.Range("L129").Value = 40427.625
$ws.Range("N129").Value = -50427.625
# Row 132
$ws.Range("H132").Value = 1061.3096
$ws.Range("I132").Value = 731.8108
$ws.Range("J132").Value = 3499.6
$ws.Range("K132").Value = 2195.4324
$ws.Range("L132").Value = 10498.8
$ws.Range("M132").Value = 334.5676000000003
$ws.Range("N132").Value = -15558.8
# Row 136
$ws.Range("H136").Value = 722.4375
$ws.Range("I136").Value = 546.7143
$ws.Range("J136").Value = 1952.5
$ws.Range("K136").Value = 1640.1429
$ws.Range("L136").Value = 5857.5
$ws.Range("M136").Value = 909.8571000000002
$ws.Range("N136").Value = -10957.5
